# Update "想去人数" (F column) counts across sheets, per the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 897
$ws1.Range("F9").Value = 1340
$ws1.Range("F14").Value = 80
$ws1.Range("F15").Value = 80
$ws1.Range("F16").Value = 1349
$ws1.Range("F27").Value = 1196

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 34

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 897
$ws4.Range("F11").Value = 1340
$ws4.Range("F19").Value = 80
$ws4.Range("F20").Value = 80
$ws4.Range("F21").Value = 1349
$ws4.Range("F32").Value = 34
$ws4.Range("F39").Value = 1196
